# ---------------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx edit replay
#
# The authoritative diff shows two logical changes:
#   1. The table on slide 16 gets a new (built-in) table style applied -
#      tableStyleId {7C54102D-3C78-4E96-A59D-3A30CBA32A11} (the deck's
#      custom "Table_0" style) -> {78339827-FF59-468C-8ABC-0F7AF0F299A4}
#      (PowerPoint's built-in "Themed Style 1 - Accent 1" gallery style).
#   2. The deck's colour theme is switched from "Integral" to the built-in
#      "Office Theme" palette (what used to live in ppt/theme/theme1.xml,
#      used only by the Notes Master, becomes the active slide-master
#      theme palette; only the 12 theme colours actually differ between
#      the two themes - fonts/format scheme are identical).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 16 table: apply the built-in "Themed Style 1 - Accent 1" style.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{78339827-FF59-468C-8ABC-0F7AF0F299A4}")
    }
}

# ---------------------------------------------------------------------------
# 2) Re-colour the presentation's theme to the "Office Theme" palette
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in that MsoThemeColorSchemeIndex
#    order). This lands on the deck's single addressable Design/theme part.
# ---------------------------------------------------------------------------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
